$wb = $excel.ActiveWorkbook

# The "Devices" sheet (first sheet) holds the cable-capacitance test data
# that needs updated expected values.
$ws = $wb.Worksheets.Item("Devices")

# Update computed AC Units values (rows 8 and 9, columns G and H)
$ws.Range("G8").Value = 20
$ws.Range("H8").Value = 172
$ws.Range("G9").Value = 56
$ws.Range("H9").Value = 153

# Move the active selection on the Devices sheet from G2 to A9
$ws.Range("A9").Select()
